$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 2).Value = 'Bachelor’s degree'
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 4
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 4
$ws.Cells.Item(3, 7).Value = 13.8

$ws.Cells.Item(4, 2).Value = 'Graduate or professional degree (MA, MS, MBA, PhD, JD, MD, DDS etc.)'
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(4, 4).Value = 10
$ws.Cells.Item(4, 5).Value = 10
$ws.Cells.Item(4, 6).Value = 24
$ws.Cells.Item(4, 7).Value = 82.8

$ws.Cells.Item(7, 2).Value = 'Academia'
$ws.Cells.Item(7, 3).Value = 4
$ws.Cells.Item(7, 4).Value = 12
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = 26
$ws.Cells.Item(7, 7).Value = 89.7

$ws.Cells.Item(8, 2).Value = 'Industry'
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(8, 7).Value = 6.9

$ws.Cells.Item(9, 2).Value = 'Postdoc'
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 2
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 4
$ws.Cells.Item(9, 7).Value = 13.8

$ws.Cells.Item(10, 2).Value = 'Staff member (including research/academic/teaching staff)'
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = 4
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 9
$ws.Cells.Item(10, 7).Value = 31

$ws.Cells.Item(11, 2).Value = 'Graduate student (including professional school student)'
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 4
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 4
$ws.Cells.Item(11, 7).Value = 13.8

$ws.Cells.Item(12, 2).Value = 'Other (please specify below)'
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 2
$ws.Cells.Item(12, 7).Value = 6.9

$ws.Cells.Item(13, 2).Value = 'Faculty member'
$ws.Cells.Item(13, 3).Value = 1
$ws.Cells.Item(13, 4).Value = 3
$ws.Cells.Item(13, 5).Value = 5
$ws.Cells.Item(13, 6).Value = 9
$ws.Cells.Item(13, 7).Value = 31

$ws.Cells.Item(14, 2).Value = 'Undergraduate student'
$ws.Cells.Item(14, 3).Value = 1
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 3.4

$ws.Cells.Item(15, 2).Value = 'Less than 1 year'
$ws.Cells.Item(15, 3).Value = 3
$ws.Cells.Item(15, 4).Value = 2
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 5
$ws.Cells.Item(15, 7).Value = 17.2

$ws.Cells.Item(16, 2).Value = '1 to 5 years'
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 4).Value = 10
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 15
$ws.Cells.Item(16, 7).Value = 51.7

$ws.Cells.Item(17, 2).Value = 'More than 5 years'
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 2
$ws.Cells.Item(17, 5).Value = 7
$ws.Cells.Item(17, 6).Value = 9
$ws.Cells.Item(17, 7).Value = 31

